$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2372.1667
$ws.Range("I40").Value = 3210.875
$ws.Range("K40").Value = 3210.875
$ws.Range("M40").Value = -3035.875

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 396383.7
$ws.Range("I64").Value = 641129.75
$ws.Range("J64").Value = 4790
$ws.Range("K64").Value = 641129.75
$ws.Range("L64").Value = 4790
$ws.Range("M64").Value = -640881.75
$ws.Range("N64").Value = -5286

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 396383.7
$ws.Range("I67").Value = 641129.75
$ws.Range("J67").Value = 4790
$ws.Range("K67").Value = 641129.75
$ws.Range("L67").Value = 4790
$ws.Range("M67").Value = -640271.75
$ws.Range("N67").Value = -6506

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4592.9443
$ws.Range("I74").Value = 4237.5
$ws.Range("J74").Value = 4877.3
$ws.Range("K74").Value = 4237.5
$ws.Range("L74").Value = 4877.3
$ws.Range("M74").Value = -3301.5
$ws.Range("N74").Value = -6749.3

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3659.4688
$ws.Range("I76").Value = 3475.92
$ws.Range("J76").Value = 4315
$ws.Range("K76").Value = 3475.92
$ws.Range("L76").Value = 4315
$ws.Range("M76").Value = -3160.92
$ws.Range("N76").Value = -4945

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4592.9443
$ws.Range("I77").Value = 4237.5
$ws.Range("J77").Value = 4877.3
$ws.Range("K77").Value = 21187.5
$ws.Range("L77").Value = 24386.5
$ws.Range("M77").Value = -16507.5
$ws.Range("N77").Value = -33746.5

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3659.4688
$ws.Range("I79").Value = 3475.92
$ws.Range("J79").Value = 4315
$ws.Range("K79").Value = 3475.92
$ws.Range("L79").Value = 4315
$ws.Range("M79").Value = -2383.92
$ws.Range("N79").Value = -6499

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4577.0386
$ws.Range("I137").Value = 2464.5
$ws.Range("J137").Value = 7041.6665
$ws.Range("K137").Value = 7393.5
$ws.Range("L137").Value = 21124.9995
$ws.Range("M137").Value = -4843.5
$ws.Range("N137").Value = -26224.9995

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3370.6562
$ws.Range("I138").Value = 1864.1482
$ws.Range("J138").Value = 3960.1594
$ws.Range("K138").Value = 5592.444600000001
$ws.Range("L138").Value = 11880.4782
$ws.Range("M138").Value = -452.4446000000007
$ws.Range("N138").Value = -22160.4782

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5894.6704
$ws.Range("I32").Value = 5285.9517
$ws.Range("J32").Value = 15999.4
$ws.Range("K32").Value = 5285.9517
$ws.Range("L32").Value = 15999.4
$ws.Range("M32").Value = -4998.9517
$ws.Range("N32").Value = -16573.4

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7724.771
$ws.Range("I61").Value = 3431.4055
$ws.Range("J61").Value = 22166.092
$ws.Range("K61").Value = 3431.4055
$ws.Range("L61").Value = 22166.092
$ws.Range("M61").Value = -3219.4055
$ws.Range("N61").Value = -22590.092

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5647.154
$ws.Range("J88").Value = 1830.1
$ws.Range("L88").Value = 1830.1
$ws.Range("N88").Value = -2642.1

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 5647.154
$ws.Range("J91").Value = 1830.1
$ws.Range("L91").Value = 1830.1
$ws.Range("N91").Value = -4638.1

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6254.528
$ws.Range("I132").Value = 2532.077
$ws.Range("J132").Value = 8358.521000000001
$ws.Range("K132").Value = 7596.231000000001
$ws.Range("L132").Value = 25075.563
$ws.Range("M132").Value = -5066.231000000001
$ws.Range("N132").Value = -30135.563

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7724.771
$ws.Range("I136").Value = 3431.4055
$ws.Range("J136").Value = 22166.092
$ws.Range("K136").Value = 10294.2165
$ws.Range("L136").Value = 66498.276
$ws.Range("M136").Value = -7744.216499999999
$ws.Range("N136").Value = -71598.276

# BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 32590.5

# BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 32590.5

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2348.9375
$ws.Range("I86").Value = 2516.3333
$ws.Range("J86").Value = 1846.75
$ws.Range("K86").Value = 2516.3333
$ws.Range("L86").Value = 1846.75
$ws.Range("M86").Value = -1393.3333
$ws.Range("N86").Value = -4092.75

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2348.9375
$ws.Range("I89").Value = 2516.3333
$ws.Range("J89").Value = 1846.75
$ws.Range("K89").Value = 12581.6665
$ws.Range("L89").Value = 9233.75
$ws.Range("M89").Value = -6965.666499999999
$ws.Range("N89").Value = -20465.75

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1776.8286
$ws.Range("I94").Value = 1695.0834
$ws.Range("K94").Value = 1695.0834
$ws.Range("M94").Value = -1244.0834

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2150.4443
$ws.Range("I99").Value = 2275.8333
$ws.Range("J99").Value = 1899.6666
$ws.Range("K99").Value = 2275.8333
$ws.Range("L99").Value = 1899.6666
$ws.Range("M99").Value = -777.8332999999998
$ws.Range("N99").Value = -4895.6666

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 42968.04
$ws.Range("I134").Value = 3124.158
$ws.Range("K134").Value = 9372.474
$ws.Range("M134").Value = -6837.474

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3072.1428
$ws.Range("I62").Value = 3000.8333
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3000.8333
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2376.8333
$ws.Range("N62").Value = -4748

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3072.1428
$ws.Range("I65").Value = 3000.8333
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 15004.1665
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -11884.1665
$ws.Range("N65").Value = -23740

# CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 10543.333
$ws.Range("I69").Value = 10543.333
$ws.Range("K69").Value = 10543.333
$ws.Range("M69").Value = -9794.333000000001

# CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 10543.333
$ws.Range("I72").Value = 10543.333
$ws.Range("K72").Value = 31629.999
$ws.Range("M72").Value = -27885.999

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1641.6461
$ws.Range("I132").Value = 1487.9395
$ws.Range("J132").Value = 1800.1562
$ws.Range("K132").Value = 4463.818499999999
$ws.Range("L132").Value = 5400.4686
$ws.Range("M132").Value = -1933.818499999999
$ws.Range("N132").Value = -10460.4686

# CUL row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3676.0605
$ws.Range("J127").Value = 3676.0605
$ws.Range("L127").Value = 11028.1815
$ws.Range("N127").Value = -20948.1815

# GSM row 52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5459.407
$ws.Range("I70").Value = 4714.143
$ws.Range("K70").Value = 4714.143
$ws.Range("M70").Value = -4444.143

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5459.407
$ws.Range("I73").Value = 4714.143
$ws.Range("K73").Value = 4714.143
$ws.Range("M73").Value = -3778.143

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11983.5
$ws.Range("I80").Value = 26500
$ws.Range("J80").Value = 4725.25
$ws.Range("K80").Value = 26500
$ws.Range("L80").Value = 4725.25
$ws.Range("M80").Value = -25502
$ws.Range("N80").Value = -6721.25

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 11983.5
$ws.Range("I83").Value = 26500
$ws.Range("J83").Value = 4725.25
$ws.Range("K83").Value = 132500
$ws.Range("L83").Value = 23626.25
$ws.Range("M83").Value = -127508
$ws.Range("N83").Value = -33610.25

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 87080.46000000001
$ws.Range("I132").Value = 127504.375
$ws.Range("J132").Value = 22402.2
$ws.Range("K132").Value = 382513.125
$ws.Range("L132").Value = 67206.60000000001
$ws.Range("M132").Value = -379983.125
$ws.Range("N132").Value = -72266.60000000001

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1373.8182
$ws.Range("I126").Value = 1234.6666
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3703.9998
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1233.9998
$ws.Range("N126").Value = -10940
